# ng: update lf forms
# Converts the Pre-TAS LF Participant form to "V2": renames the EU field,
# relabels the cluster fields from "school" to "community", bumps the
# repeat/form identifiers to the _v2 variants, trims the "team" choice
# list down to 01-04/99, and updates the form title/id on the settings
# sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "survey" sheet
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 3: p_eu -> p_region (label "Select your EU" is unchanged)
$survey.Cells.Item(3, 2).Value = "p_region"

# Row 5: p_cluster_name label "Select the school" -> "Select the community"
$survey.Cells.Item(5, 3).Value = "Select the community"

# Row 6: p_cluster_id label "Enter the school ID" -> "Enter the community ID"
$survey.Cells.Item(6, 3).Value = "Enter the community ID"

# Row 8: repeat name ng_pretas_p_202404 -> ng_pretas_p_202404_v2
$survey.Cells.Item(8, 2).Value = "ng_pretas_p_202404_v2"

# Move the saved selection to B3 (matches the author's recorded cursor pos)
$survey.Range("B3").Select()

# ---------------------------------------------------------------------
# "choices" sheet - shrink the "team" list from 01..10,99 down to 01..04,99
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Rows("6:11").Delete()
$choices.Range("A6:XFD11").Select()

# ---------------------------------------------------------------------
# "settings" sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Cells.Item(2, 1).Value = "(Apr 2024) - 2. Nigeria - Pre TAS LF Participant Form V2"
$settings.Cells.Item(2, 2).Value = "ng_lf_pretas_202404_2_participant_v2"
$settings.Range("B2").Select()

$survey.Activate()
